# Add a second responsible professor ("Docentes responsáveis:") to the
# LOM3210 course sheet: insert a new row right below the existing
# "519033 - Carlos Yujiro Shigue" entry (row 13) and fill it with the new
# professor's name in columns B and C, matching the layout of every other
# two-column (current/modified) content row on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push row 14 (and everything below it) down by one row.
$ws.Rows("14:14").Insert()

# Populate the newly inserted row with the new professor.
$ws.Range("B14").Value = "1176388 - Luiz Tadeu Fernandes Eleno"
$ws.Range("C14").Value = "1176388 - Luiz Tadeu Fernandes Eleno"
